# Auto-generated edit script: refresh cryptos list snapshot (prices / 1h volume %).
# Rows 20-51 shift down by one (a new "BitDAO" entry was inserted after Avalanche),
# and every Price / Volume(1h) cell is refreshed to the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '31.156.42'
$ws.Range('E2').Value = '  +2.18%  '

# Row 3
$ws.Range('D3').Value = '1.970.30'
$ws.Range('E3').Value = '  +3.20%  '

# Row 4
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '1.006'
$c.Style = "Normal"
$ws.Range('E4').Value = '  +0.53%  '

# Row 5
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '248.51'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +1.89%  '

# Row 6
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '1.005'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +0.39%  '

# Row 7
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.4887'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +1.09%  '

# Row 8
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '44.95'
$c.Style = "Normal"
$ws.Range('E8').Value = '  +1.21%  '

# Row 9
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.2954'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +2.49%  '

# Row 10
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.06839'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +0.65%  '

# Row 11
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '19.20'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -0.64%  '

# Row 12
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '107.55'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -3.52%  '

# Row 13
$ws.Range('D13').Value = '1.966.66'
$ws.Range('E13').Value = '  +2.92%  '

# Row 14
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '0.07799'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +3.12%  '

# Row 15
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '5.444'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +1.09%  '

# Row 16
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '0.7066'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +5.84%  '

# Row 17
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '285.62'
$c.Style = "Normal"

# Row 18
$ws.Range('D18').Value = '31.154.94'
$ws.Range('E18').Value = '  +2.14%  '

# Row 19
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '13.31'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +2.56%  '

# Row 20
$ws.Range('B20').Value = 'BitDAO'
$ws.Range('C20').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '0.4971'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +14.54%  '

# Row 21
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '0.000007751'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +2.52%  '

# Row 22
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.219.80'
$ws.Range('E22').Value = '  +2.66%  '

# Row 23
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '5.622'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +2.47%  '

# Row 24
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E24').Value = '  +0.36%  '

# Row 25
$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '1.007'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +0.78%  '

# Row 26
$ws.Range('B26').Value = 'Chainlink'
$ws.Range('C26').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '6.642'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +4.05%  '

# Row 27
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '10.00'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +5.95%  '

# Row 28
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '170.36'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +3.07%  '

# Row 29
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '20.13'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -0.34%  '

# Row 30
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '2.191'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +6.17%  '

# Row 31
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '0.1067'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +0.21%  '

# Row 32
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '1.443'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +1.05%  '

# Row 33
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '4.835'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +19.49%  '

# Row 34
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '4.518'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +9.59%  '

# Row 35
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.05091'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +2.46%  '

# Row 36
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.7708'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +4.99%  '

# Row 37
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '1.172'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +3.79%  '

# Row 38
$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '2.744'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +1.33%  '

# Row 39
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.02044'
$c.Style = "Normal"
$ws.Range('E39').Value = '  +0.74%  '

# Row 40
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '2.735'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +2.07%  '

# Row 41
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '6.469'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +12.06%  '

# Row 42
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '2.129'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +6.10%  '

# Row 43
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.8891'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +3.26%  '

# Row 44
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '110.10'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +0.77%  '

# Row 45
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.4473'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +1.23%  '

# Row 46
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '73.38'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +6.42%  '

# Row 47
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +0.37%  '

# Row 48
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '7.529'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +4.97%  '

# Row 49
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '992.28'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +18.02%  '

# Row 50
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '9.419'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +2.85%  '

# Row 51
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.1270'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +3.67%  '

